# Testing Matrix.xlsx
# - split the "MinMax " column into "MinMax 1" / "MinMax 2" in both tables
# - shift "Sigma" one column to the right to make room, and extend the
#   green-checkmark result blocks to cover it (plus one spacer column)
# - mark every existing result cell with a passing checkmark
# - move the legend (check/cross) two columns over, and add a note about
#   the new two-result test

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$check = "✔️"
$cross = "✘"

$xlCenter = -4108
$xlPasteFormats = -4122
$greenFont = 5287936   # RGB(0,176,80) == style used for check marks
$redFont = 255         # RGB(255,0,0)  == style used for the cross

# --- Move the legend two columns to the right (H1:H2 -> J1:J2), keeping
#     its existing (green / red) formatting ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("J1").Value = $check
$ws.Range("J2").Value = $cross
$ws.Range("H1:H2").Clear()

# --- New note next to the first result block (added early so the shared
#     string table picks up the same ordering as the authored workbook) ---
$ws.Range("I4").Value = "Two tests: complet and cancel"

# --- Header rows (3 and 11): "MinMax " -> "MinMax 1" / "MinMax 2", and
#     "Sigma" shifts one column right, from E to F ---
foreach ($r in 3, 11) {
    $ws.Range("F${r}").Value = "Sigma"
    $ws.Range("F${r}").HorizontalAlignment = $xlCenter

    $ws.Range("D${r}").Value = "MinMax 1"
    $ws.Range("E${r}").Value = "MinMax 2"
}

# --- Result blocks: fill B:F with checkmarks and add a blank styled
#     spacer cell in G, for each of the existing pass/fail rows ---
foreach ($r in 4, 5, 6, 7, 12, 13) {
    $rng = $ws.Range("B${r}:F${r}")
    $rng.Value = $check
    $rng.HorizontalAlignment = $xlCenter
    $rng.Font.Color = $greenFont

    $spacer = $ws.Range("G${r}")
    $spacer.HorizontalAlignment = $xlCenter
    $spacer.Font.Color = $greenFont
}

# Row 14 only gains checkmarks in B:E (no F/G formatting for this row)
$rng14 = $ws.Range("B14:E14")
$rng14.Value = $check
$rng14.HorizontalAlignment = $xlCenter
$rng14.Font.Color = $greenFont

# Row 15 keeps its blank styled cells, extended to cover F:G as well
$rng15 = $ws.Range("F15:G15")
$rng15.HorizontalAlignment = $xlCenter
$rng15.Font.Color = $greenFont

# --- Columns: match column F's width to the existing B:E block, and add
#     narrow spacer columns (G, H:I) ---
$ws.Columns("F:F").ColumnWidth = $ws.Columns("B:B").ColumnWidth
$ws.Columns("G:I").ColumnWidth = 4.33

# --- Final selection matches the author's last edit position ---
$ws.Range("E14").Select()
